$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GUION")
Write-Host $ws.Name
